# Add team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the existing header row (bold, centered, bordered).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting (not the value) from the neighboring header cell so the
# new headers share the same style index as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Determine the last data row (mirrors the existing used range, row 66).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 77   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 85   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
